$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rln1"
$ws.Range("C2").Value = "Rxfp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.09405766666666666
$ws.Range("H2").Value = 0.282173
$ws.Range("I2").Value = 0.203423347045177
$ws.Range("J2").Value = 0.203423347045177
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0002143333333333333
$ws.Range("N2").Value = 0.000643
$ws.Range("O2").Value = 0.0005591790590486129
$ws.Range("P2").Value = 0.0005591790590486129
$ws.Range("Q2").Value = 0.00002015969322222222
$ws.Range("R2").Value = 0.000181437239
$ws.Range("S2").Value = 0.0001137500757892415
$ws.Range("T2").Value = 0.0001137500757892415

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rln1"
$ws.Range("C3").Value = "Rxfp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.09405766666666666
$ws.Range("H3").Value = 0.282173
$ws.Range("I3").Value = 0.203423347045177
$ws.Range("J3").Value = 0.203423347045177
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.380202
$ws.Range("N3").Value = 1.140606
$ws.Range("O3").Value = 0.991917558048526
$ws.Range("P3").Value = 0.991917558048526
$ws.Range("Q3").Value = 0.035760912982
$ws.Range("R3").Value = 0.321848216838
$ws.Range("S3").Value = 0.2017791896511098
$ws.Range("T3").Value = 0.2017791896511098

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rln1"
$ws.Range("C4").Value = "Rxfp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.09405766666666666
$ws.Range("H4").Value = 0.282173
$ws.Range("I4").Value = 0.203423347045177
$ws.Range("J4").Value = 0.203423347045177
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.002883666666666667
$ws.Range("N4").Value = 0.008651
$ws.Range("O4").Value = 0.007523262892425429
$ws.Range("P4").Value = 0.00752326289242543
$ws.Range("Q4").Value = 0.0002712309581111111
$ws.Range("R4").Value = 0.002441078623
$ws.Range("S4").Value = 0.00153040731827796
$ws.Range("T4").Value = 0.001530407318277961

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rln1"
$ws.Range("C5").Value = "Rxfp2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2889726666666667
$ws.Range("H5").Value = 0.866918
$ws.Range("I5").Value = 0.624976029505696
$ws.Range("J5").Value = 0.6249760295056961
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.0002143333333333333
$ws.Range("N5").Value = 0.000643
$ws.Range("O5").Value = 0.0005591790590486129
$ws.Range("P5").Value = 0.0005591790590486129
$ws.Range("Q5").Value = 0.00006193647488888889
$ws.Range("R5").Value = 0.000557428274
$ws.Range("S5").Value = 0.0003494735081069332
$ws.Range("T5").Value = 0.0003494735081069333

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rln1"
$ws.Range("C6").Value = "Rxfp2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2889726666666667
$ws.Range("H6").Value = 0.866918
$ws.Range("I6").Value = 0.624976029505696
$ws.Range("J6").Value = 0.6249760295056961
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.380202
$ws.Range("N6").Value = 1.140606
$ws.Range("O6").Value = 0.991917558048526
$ws.Range("P6").Value = 0.991917558048526
$ws.Range("Q6").Value = 0.109867985812
$ws.Range("R6").Value = 0.9888118723079999
$ws.Range("S6").Value = 0.6199246970261535
$ws.Range("T6").Value = 0.6199246970261536

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rln1"
$ws.Range("C7").Value = "Rxfp2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2889726666666667
$ws.Range("H7").Value = 0.866918
$ws.Range("I7").Value = 0.624976029505696
$ws.Range("J7").Value = 0.6249760295056961
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.002883666666666667
$ws.Range("N7").Value = 0.008651
$ws.Range("O7").Value = 0.007523262892425429
$ws.Range("P7").Value = 0.00752326289242543
$ws.Range("Q7").Value = 0.0008333008464444444
$ws.Range("R7").Value = 0.007499707618000001
$ws.Range("S7").Value = 0.004701858971435583
$ws.Range("T7").Value = 0.004701858971435584

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Rln1"
$ws.Range("C8").Value = "Rxfp2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.07934366666666666
$ws.Range("H8").Value = 0.238031
$ws.Range("I8").Value = 0.171600623449127
$ws.Range("J8").Value = 0.1716006234491271
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.0002143333333333333
$ws.Range("N8").Value = 0.000643
$ws.Range("O8").Value = 0.0005591790590486129
$ws.Range("P8").Value = 0.0005591790590486129
$ws.Range("Q8").Value = 0.00001700599255555556
$ws.Range("R8").Value = 0.000153053933
$ws.Range("S8").Value = 0.00009595547515243820
$ws.Range("T8").Value = 0.00009595547515243823

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Rln1"
$ws.Range("C9").Value = "Rxfp2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.07934366666666666
$ws.Range("H9").Value = 0.238031
$ws.Range("I9").Value = 0.171600623449127
$ws.Range("J9").Value = 0.1716006234491271
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.380202
$ws.Range("N9").Value = 1.140606
$ws.Range("O9").Value = 0.991917558048526
$ws.Range("P9").Value = 0.991917558048526
$ws.Range("Q9").Value = 0.030166620754
$ws.Range("R9").Value = 0.271499586786
$ws.Range("S9").Value = 0.1702136713712627
$ws.Range("T9").Value = 0.1702136713712628

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Rln1"
$ws.Range("C10").Value = "Rxfp2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.07934366666666666
$ws.Range("H10").Value = 0.238031
$ws.Range("I10").Value = 0.171600623449127
$ws.Range("J10").Value = 0.1716006234491271
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.002883666666666667
$ws.Range("N10").Value = 0.008651
$ws.Range("O10").Value = 0.007523262892425429
$ws.Range("P10").Value = 0.00752326289242543
$ws.Range("Q10").Value = 0.0002288006867777778
$ws.Range("R10").Value = 0.002059206181
$ws.Range("S10").Value = 0.001290996602711886
$ws.Range("T10").Value = 0.001290996602711887
